$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Protect column D (Price) as text before writing, since several values
# look like plain numbers (e.g. "571.51") and Excel would otherwise coerce
# them to the Number type, losing the original text representation
# (trailing zeros, etc.). Dotted-thousands values (e.g. "61.363.09") are
# already non-numeric and unaffected either way.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "61.363.09"
$ws.Range("E2").Value = "  +0.29%  "

$ws.Range("D3").Value = "3.372.93"
$ws.Range("E3").Value = "  +0.01%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "571.51"
$ws.Range("E5").Value = "  -0.16%  "

$ws.Range("D6").Value = "136.85"
$ws.Range("E6").Value = "  -0.39%  "

$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("D8").Value = "3.372.23"
$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("D9").Value = "0.474"
$ws.Range("E9").Value = "  -0.68%  "

$ws.Range("D10").Value = "7.47"
$ws.Range("E10").Value = "  -1.83%  "

$ws.Range("E11").Value = "  +0.30%  "

$ws.Range("D12").Value = "0.390"
$ws.Range("E12").Value = "  -0.93%  "

$ws.Range("D13").Value = "3.948.29"
$ws.Range("E13").Value = "  -0.08%  "

$ws.Range("E14").Value = "  +1.96%  "

$ws.Range("E15").Value = "  +0.49%  "

$ws.Range("D16").Value = "25.86"
$ws.Range("E16").Value = "  +2.72%  "

$ws.Range("D17").Value = "3.374.72"
$ws.Range("E17").Value = "  +0.09%  "

$ws.Range("D18").Value = "61.454.65"
$ws.Range("E18").Value = "  +0.39%  "

$ws.Range("B19").Value = "Polkadot"
$ws.Range("C19").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D19").Value = "5.88"
$ws.Range("E19").Value = "  +0.11%  "

$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").Value = "13.93"
$ws.Range("E20").Value = "  +0.11%  "

$ws.Range("E21").Value = "  -0.75%  "

$ws.Range("D22").Value = "376.18"
$ws.Range("E22").Value = "  -0.82%  "

$ws.Range("D23").Value = "0.554"
$ws.Range("E23").Value = "  -2.68%  "

$ws.Range("D24").Value = "3.512.63"
$ws.Range("E24").Value = "  +0.10%  "

$ws.Range("D25").Value = "0.999"
$ws.Range("E25").Value = "  -0.18%  "

$ws.Range("D26").Value = "71.59"
$ws.Range("E26").Value = "  +1.28%  "

$ws.Range("E27").Value = "  +4.99%  "

$ws.Range("D28").Value = "1.72"

$ws.Range("D29").Value = "7.50"
$ws.Range("E29").Value = "  -3.48%  "

$ws.Range("E30").Value = "  +0.20%  "

$ws.Range("D31").Value = "0.161"
$ws.Range("E31").Value = "  +3.38%  "

$ws.Range("D32").Value = "8.22"
$ws.Range("E32").Value = "  +0.35%  "

$ws.Range("E33").Value = "  +1.26%  "

$ws.Range("E34").Value = "  +0.06%  "

$ws.Range("D35").Value = "23.51"
$ws.Range("E35").Value = "  +0.56%  "

$ws.Range("E36").Value = "  -7.25%  "

$ws.Range("D37").Value = "6.81"
$ws.Range("E37").Value = "  -3.39%  "

$ws.Range("E38").Value = "  -1.54%  "

$ws.Range("D39").Value = "164.95"
$ws.Range("E39").Value = "  +1.59%  "

$ws.Range("D40").Value = "0.0770"
$ws.Range("E40").Value = "  -3.20%  "

$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  -0.06%  "

$ws.Range("D42").Value = "1.72"
$ws.Range("E42").Value = "  +0.30%  "

$ws.Range("D43").Value = "0.772"
$ws.Range("E43").Value = "  +1.81%  "

$ws.Range("D44").Value = "41.47"
$ws.Range("E44").Value = "  -0.13%  "

$ws.Range("D45").Value = "1.20"
$ws.Range("E45").Value = "  +0.30%  "

$ws.Range("E46").Value = "  -1.26%  "

$ws.Range("D47").Value = "24.63"
$ws.Range("E47").Value = "  +5.49%  "

$ws.Range("D48").Value = "6.83"
$ws.Range("E48").Value = "  -1.89%  "

$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").Value = "22.63"
$ws.Range("E49").Value = "  -1.66%  "

$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").Value = "2.357.06"
$ws.Range("E50").Value = "  +0.56%  "

$ws.Range("D51").Value = "0.0263"
$ws.Range("E51").Value = "  +0.02%  "

# Restore the default cell style now that the text values are committed,
# so no lingering explicit text-format style is left on the cells.
$ws.Range("D2:D51").Style = "Normal"